# This script applies the "final commit of first project" update:
#  - TaskList sheet: fixes/re-adds rows 122-126 (normalizing format) and appends new rows 127-135
#  - error report sheet: turns the informal note in (old) B22 into a full Lopa_ERR21 error-report row
#  - removes the now-unused "javax.el.ELException..." shared string / Arial Unicode MS font+style

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TaskList")
$errws = $wb.Worksheets.Item("error report")

# ---- TaskList rows 122-135 ----
# columns: A S.No | B Date | C Task | D Link to video | E Link to material | F Time Taken | G Errors Y/N | H Error/Solution ref

$rows = @(
    @{ Row=122; No=121; Date="3/12/2017"; Task="Cart WebFlow Testing";              Time="5.5Hrs" },
    @{ Row=123; No=122; Date="3/12/2017"; Task="Project2 Discussion";               Time="1hr" },
    @{ Row=124; No=123; Date="3/12/2017"; Task="Configuration for Project2";        Time="1hr" },
    @{ Row=125; No=124; Date="3/12/2017"; Task="Displaying Hello world";            Time="1 hr" },
    @{ Row=126; No=125; Date="3/12/2017"; Task="Installation of Visual Code";       Time="1hr" },
    @{ Row=127; No=126; Date="3/12/2017"; Task="Testing Angular Application";       Time="1hr" },
    @{ Row=128; No=127; Date="3/13/2017"; Task="Custom exception handler example";  Time="30Minutes" },
    @{ Row=129; No=128; Date="3/13/2017"; Task="Linking angular to spring app";     Time="3hrs" },
    @{ Row=130; No=129; Date="3/13/2017"; Task="Error Solving";                     Time="60Minutes" },
    @{ Row=131; No=130; Date="3/13/2017"; Task="Order Invoice page";                Time="60Minutes" },
    @{ Row=132; No=131; Date="3/13/2017"; Task="Webflow checkout";                  Time="120Minutes" },
    @{ Row=133; No=132; Date="3/13/2017"; Task="Update Product Item";               Time="60Minutes" },
    @{ Row=134; No=133; Date="3/14/2017"; Task="Cart Error Resolve";                Time="30Minutes" },
    @{ Row=135; No=134; Date="3/14/2017"; Task="Update Product Item";               Time="30Minutes"; Err="Yes"; ErrId="Lopa_ERR21" }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value2 = $r.No
    $ws.Cells.Item($row, 2).Value2 = $r.Date
    $ws.Cells.Item($row, 3).Value2 = $r.Task
    $ws.Cells.Item($row, 6).Value2 = $r.Time
    if ($r.ContainsKey("Err")) {
        $ws.Cells.Item($row, 7).Value2 = $r.Err
        $ws.Cells.Item($row, 8).Value2 = $r.ErrId
    }
}

# ---- error report row 22: replace informal note with a full Lopa_ERR21 entry ----
$errws.Cells.Item(22, 1).Value2 = "Lopa_ERR21"
$errws.Cells.Item(22, 2).Value2 = "Statestate Exception in updating product item"
$errws.Cells.Item(22, 3).Value2 = 'Add the following code in side the productCRUD.jsp page within form tag  <form:hidden path="productId">'
